# Update "想去人数" (want-to-go count) figures in column F for the
# sheets that list individual events ("展览" and "全部类型").
# Both sheets mirror the same underlying data, so the same F-column
# cells need the same updates applied on each of them.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F7"  = 2083
    "F8"  = 73
    "F10" = 4554
    "F13" = 109
    "F19" = 3436
    "F21" = 549
    "F31" = 681
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
